# Update gh-pages to output generated at 456a3b4
# This script updates the "想去人数" (F column) values across the four
# worksheets of the workbook to match the freshly regenerated data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 1746
$ws.Range("F5").Value  = 437
$ws.Range("F7").Value  = 60
$ws.Range("F8").Value  = 323
$ws.Range("F9").Value  = 288
$ws.Range("F10").Value = 1673
$ws.Range("F11").Value = 331
$ws.Range("F12").Value = 1382
$ws.Range("F13").Value = 783
$ws.Range("F14").Value = 316
$ws.Range("F15").Value = 652
$ws.Range("F16").Value = 12595
$ws.Range("F17").Value = 12629
$ws.Range("F18").Value = 932
$ws.Range("F19").Value = 730
$ws.Range("F21").Value = 295
$ws.Range("F23").Value = 484
$ws.Range("F24").Value = 1974
$ws.Range("F25").Value = 21
$ws.Range("F27").Value = 225
$ws.Range("F28").Value = 657

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value  = 69
$ws.Range("F6").Value  = 9
$ws.Range("F8").Value  = 125
$ws.Range("F9").Value  = 48
$ws.Range("F10").Value = 65

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 147

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 147
$ws.Range("F6").Value  = 1746
$ws.Range("F7").Value  = 437
$ws.Range("F10").Value = 60
$ws.Range("F12").Value = 323
$ws.Range("F14").Value = 288
$ws.Range("F15").Value = 1673
$ws.Range("F16").Value = 331
$ws.Range("F17").Value = 1382
$ws.Range("F18").Value = 783
$ws.Range("F19").Value = 316
$ws.Range("F20").Value = 69
$ws.Range("F21").Value = 652
$ws.Range("F22").Value = 12595
$ws.Range("F23").Value = 12629
$ws.Range("F24").Value = 932
$ws.Range("F25").Value = 730
$ws.Range("F27").Value = 295
$ws.Range("F29").Value = 484
$ws.Range("F30").Value = 9
$ws.Range("F32").Value = 1974
$ws.Range("F33").Value = 21
$ws.Range("F34").Value = 125
$ws.Range("F36").Value = 48
$ws.Range("F37").Value = 225
$ws.Range("F38").Value = 657
$ws.Range("F39").Value = 65

$wb.Save()
